$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: add new columns P1 and Q1, styled like other header cells (s="1")
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Rows 2-25: swap values in columns I, K, M, O (1<->2), and add columns P and Q with value 2
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value()   # column I
    $kVal = $ws.Cells.Item($r, 11).Value()  # column K
    $mVal = $ws.Cells.Item($r, 13).Value()  # column M
    $oVal = $ws.Cells.Item($r, 15).Value()  # column O

    if ($iVal -eq 1) { $ws.Cells.Item($r, 9).Value = 2 } else { $ws.Cells.Item($r, 9).Value = 1 }
    if ($kVal -eq 1) { $ws.Cells.Item($r, 11).Value = 2 } else { $ws.Cells.Item($r, 11).Value = 1 }
    if ($mVal -eq 1) { $ws.Cells.Item($r, 13).Value = 2 } else { $ws.Cells.Item($r, 13).Value = 1 }
    if ($oVal -eq 1) { $ws.Cells.Item($r, 15).Value = 2 } else { $ws.Cells.Item($r, 15).Value = 1 }

    $ws.Cells.Item($r, 16).Value = 2  # column P
    $ws.Cells.Item($r, 17).Value = 2  # column Q
}
